$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Dereck Lively II", "C", "Dallas Mavericks"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Norman Powell", "SG,SF", "LA Clippers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
